$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 601
$ws.Range("I29").Value = 450
$ws.Range("J29").Value = 752
$ws.Range("K29").Value = 1350
$ws.Range("L29").Value = 2256
$ws.Range("M29").Value = -1069
$ws.Range("N29").Value = -2818

$ws.Range("H38").Value = 101638.1
$ws.Range("I38").Value = 375171.62
$ws.Range("J38").Value = 2171.3635
$ws.Range("K38").Value = 1125514.86
$ws.Range("L38").Value = 6514.0905
$ws.Range("M38").Value = -1125142.86
$ws.Range("N38").Value = -7258.0905

$ws.Range("H40").Value = 1453
$ws.Range("I40").Value = 1354.4546
$ws.Range("J40").Value = 1633.6666
$ws.Range("K40").Value = 1354.4546
$ws.Range("L40").Value = 1633.6666
$ws.Range("M40").Value = -1179.4546
$ws.Range("N40").Value = -1983.6666

$ws.Range("H53").Value = 6944864.5
$ws.Range("I53").Value = 13889110
$ws.Range("J53").Value = 618.875
$ws.Range("K53").Value = 13889110
$ws.Range("L53").Value = 618.875
$ws.Range("M53").Value = -13888473
$ws.Range("N53").Value = -1892.875

$ws.Range("H58").Value = 817.4666999999999
$ws.Range("I58").Value = 65.833336
$ws.Range("J58").Value = 1318.5555
$ws.Range("K58").Value = 197.500008
$ws.Range("L58").Value = 3955.6665
$ws.Range("M58").Value = -47.50000800000001
$ws.Range("N58").Value = -4255.666499999999

$ws.Range("H74").Value = 8563
$ws.Range("I74").Value = 19250
$ws.Range("K74").Value = 19250
$ws.Range("M74").Value = -18314

$ws.Range("H77").Value = 8563
$ws.Range("I77").Value = 19250
$ws.Range("K77").Value = 96250
$ws.Range("M77").Value = -91570

$ws.Range("H87").Value = 11981.675
$ws.Range("J87").Value = 11981.675
$ws.Range("L87").Value = 11981.675
$ws.Range("N87").Value = -14477.675

$ws.Range("H90").Value = 11981.675
$ws.Range("J90").Value = 11981.675
$ws.Range("L90").Value = 35945.02499999999
$ws.Range("N90").Value = -48425.02499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8335.333000000001
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 9402.4
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 9402.4
$ws.Range("M102").Value = -1378
$ws.Range("N102").Value = -12646.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29443060
$ws.Range("I86").Value = 38463396
$ws.Range("J86").Value = 126961.375
$ws.Range("K86").Value = 38463396
$ws.Range("L86").Value = 126961.375
$ws.Range("M86").Value = -38462273
$ws.Range("N86").Value = -129207.375

$ws.Range("H89").Value = 29443060
$ws.Range("I89").Value = 38463396
$ws.Range("J89").Value = 126961.375
$ws.Range("K89").Value = 192316980
$ws.Range("L89").Value = 634806.875
$ws.Range("M89").Value = -192311364
$ws.Range("N89").Value = -646038.875

$ws.Range("H105").Value = 4174.6924
$ws.Range("I105").Value = 3160
$ws.Range("K105").Value = 3160
$ws.Range("M105").Value = -1413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 400
$ws.Range("J17").Value = 550
$ws.Range("L17").Value = 1650
$ws.Range("N17").Value = -1988

$ws.Range("H114").Value = 561754.5600000001
$ws.Range("I114").Value = 7984
$ws.Range("J114").Value = 2001558
$ws.Range("K114").Value = 23952
$ws.Range("L114").Value = 6004674
$ws.Range("M114").Value = -20698
$ws.Range("N114").Value = -6011182

$ws.Range("H121").Value = 66028.67999999999
$ws.Range("I121").Value = 5298.3335
$ws.Range("J121").Value = 80603.96000000001
$ws.Range("K121").Value = 15895.0005
$ws.Range("L121").Value = 241811.88
$ws.Range("M121").Value = -14585.0005
$ws.Range("N121").Value = -244431.88

$ws.Range("H122").Value = 1028.7778
$ws.Range("J122").Value = 1351.5
$ws.Range("L122").Value = 12163.5
$ws.Range("N122").Value = -17063.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 7800
$ws.Range("J62").Value = 7800
$ws.Range("L62").Value = 7800
$ws.Range("N62").Value = -9172

$ws.Range("H65").Value = 7800
$ws.Range("J65").Value = 7800
$ws.Range("L65").Value = 23400
$ws.Range("N65").Value = -30264

$ws.Range("H132").Value = 50066.383
$ws.Range("I132").Value = 2048.7856
$ws.Range("J132").Value = 146101.58
$ws.Range("K132").Value = 6146.3568
$ws.Range("L132").Value = 438304.74
$ws.Range("M132").Value = -3616.3568
$ws.Range("N132").Value = -443364.74

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1670.659
$ws.Range("I40").Value = 1621.1177
$ws.Range("J40").Value = 1839.1
$ws.Range("K40").Value = 1621.1177
$ws.Range("L40").Value = 1839.1
$ws.Range("M40").Value = -1485.1177
$ws.Range("N40").Value = -2111.1

$ws.Range("H68").Value = 1937.2106
$ws.Range("I68").Value = 1831.0769
$ws.Range("J68").Value = 2167.1667
$ws.Range("K68").Value = 1831.0769
$ws.Range("L68").Value = 2167.1667
$ws.Range("M68").Value = -1082.0769
$ws.Range("N68").Value = -3665.1667

$ws.Range("H71").Value = 1937.2106
$ws.Range("I71").Value = 1831.0769
$ws.Range("J71").Value = 2167.1667
$ws.Range("K71").Value = 9155.3845
$ws.Range("L71").Value = 10835.8335
$ws.Range("M71").Value = -5411.3845
$ws.Range("N71").Value = -18323.8335

$ws.Range("H100").Value = 2542
$ws.Range("I100").Value = 1083.3334
$ws.Range("J100").Value = 4000.6667
$ws.Range("K100").Value = 1083.3334
$ws.Range("L100").Value = 4000.6667
$ws.Range("M100").Value = -542.3334
$ws.Range("N100").Value = -5082.6667

$ws.Range("H132").Value = 670203.4
$ws.Range("I132").Value = 1669224.6
$ws.Range("J132").Value = 4189.222
$ws.Range("K132").Value = 5007673.800000001
$ws.Range("L132").Value = 12567.666
$ws.Range("M132").Value = -5005143.800000001
$ws.Range("N132").Value = -17627.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 46653.332
$ws.Range("J131").Value = 46653.332
$ws.Range("L131").Value = 46653.332
$ws.Range("N131").Value = -56733.332

$ws.Range("H132").Value = 89196024
$ws.Range("I132").Value = 281251840
$ws.Range("J132").Value = 3837875.5
$ws.Range("K132").Value = 843755520
$ws.Range("L132").Value = 11513626.5
$ws.Range("M132").Value = -843752990
$ws.Range("N132").Value = -11518686.5
